$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.608.12"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.478.21"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.42"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.15"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.479.00"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +5.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.31"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.080.42"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.64"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.600.03"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.472.91"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.28"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.40"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.25"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.71"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.44"
$ws.Range("E30").Value = "  +7.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.70"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.12"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.88"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.97"
$ws.Range("E39").Value = "  +5.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.050.99"
$ws.Range("E40").Value = "  +5.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0773"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.20"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.51"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.84"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.775"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.03"
$ws.Range("E47").Value = "  +10.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.72"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "311.40"
$ws.Range("E51").Value = "  +4.44%  "
